$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.281.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "'3.541.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'603.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("D6").Value = "'140.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.56%  "
$ws.Range("D7").Value = "'3.541.54"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.70%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").Value = "'0.126"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("E11").Value = "  -5.06%  "
$ws.Range("E12").Value = "  +3.66%  "
$ws.Range("D13").Value = "'4.137.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").Value = "'0.0000188"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.96%  "
$ws.Range("D15").Value = "'27.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "'3.553.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "'65.379.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'10.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.51%  "
$ws.Range("D20").Value = "'5.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.70%  "
$ws.Range("D21").Value = "'14.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.32%  "
$ws.Range("D22").Value = "'396.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").Value = "'0.574"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.80%  "
$ws.Range("D24").Value = "'3.683.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("D25").Value = "'74.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +9.86%  "
$ws.Range("D28").Value = "'7.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.71%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'2.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.93%  "
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").Value = "'3.548.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.69%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "'23.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("D36").Value = "'1.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.23%  "
$ws.Range("D37").Value = "'7.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("D38").Value = "'170.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("D39").Value = "'1.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.38%  "
$ws.Range("E40").Value = "  +3.26%  "
$ws.Range("D41").Value = "'0.0817"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.95%  "
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("E43").Value = "  +16.40%  "
$ws.Range("D44").Value = "'42.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.80%  "
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "'4.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("E47").Value = "  +9.62%  "
$ws.Range("D49").Value = "'2.454.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.40%  "
$ws.Range("D50").Value = "'6.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.73%  "
$ws.Range("D51").Value = "'2.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +15.84%  "
